# Financiero.xlsx edit script
# Summary of the change (per commit "Mejora en deseño y visualización"):
#   - Rename the "Ahorro" concept to "Abono" throughout (column values and
#     the "Total Ahorro" headers / table column on the summary sheets).
#   - Remove the trailing test rows ("Pepito" / "Prueba 3") that were added
#     to the Proveedores sheet and that leaked into the Resumen summary.
#   - Give row 2 of Proveedores (the very first date) a date-only display
#     format, while the remaining dates keep their date+time format.
#   - Refresh the sequential "Id" numbering on the Resumen sheet now that
#     rows were removed.
#   - Update the active sheet/selection bookkeeping so ResumenCliente ends
#     up the selected tab.

$wb = $excel.ActiveWorkbook

$wsProveedores       = $wb.Worksheets.Item("Proveedores")
$wsResumen           = $wb.Worksheets.Item("Resumen")
$wsProveedoresCliente = $wb.Worksheets.Item("ProveedoresCliente")
$wsResumenCliente    = $wb.Worksheets.Item("ResumenCliente")

# ---------------------------------------------------------------------
# 1) Proveedores: drop the trailing test rows (Pepito / Prueba 3), 148
#    data rows is the real data set (rows 2-148).
# ---------------------------------------------------------------------
$wsProveedores.Range("A149:A154").EntireRow.Delete() | Out-Null

# Rename "Ahorro" -> "Abono" in the Detalle column (column D).
$wsProveedores.Range("D2:D148").Replace("Ahorro", "Abono") | Out-Null

# Row 2's date should render as a plain date, the rest keep date+time.
$wsProveedores.Range("B2").NumberFormat = "yyyy-mm-dd"
$wsProveedores.Range("B3:B148").NumberFormat = "yyyy-mm-dd h:mm:ss"

$wsProveedores.Select() | Out-Null
$wsProveedores.Range("D13").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) Resumen: rename the header/table column, drop the two stray rows
#    that came from the removed test data, and renumber the Id column.
# ---------------------------------------------------------------------
$wsResumen.Range("D1").Value = "Total Abonos"

$wsResumen.Range("A2").Value = 1
$wsResumen.Range("A3").Value = 2
$wsResumen.Range("A4").Value = 3
$wsResumen.Range("A5").Value = 4
$wsResumen.Range("A6").Value = 5
$wsResumen.Range("A7").Value = 6

$wsResumen.Range("A8:A9").EntireRow.Delete() | Out-Null

$wsResumen.ListObjects.Item(1).ListColumns.Item("Total Ahorro").Name = "Total Abonos"

$wsResumen.Select() | Out-Null
$wsResumen.Range("H5").Select() | Out-Null

# ---------------------------------------------------------------------
# 3) ProveedoresCliente: just a selection bookkeeping change.
# ---------------------------------------------------------------------
$wsProveedoresCliente.Select() | Out-Null
$wsProveedoresCliente.Range("D11").Select() | Out-Null

# ---------------------------------------------------------------------
# 4) ResumenCliente: rename header, becomes the active/selected sheet.
# ---------------------------------------------------------------------
$wsResumenCliente.Range("D1").Value = "Total Abonos"

$wsResumenCliente.Select() | Out-Null
$wsResumenCliente.Range("D2").Select() | Out-Null
